$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '29.543.74'
Set-TextValue $ws.Range('E2') '  +0.22%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.913.59'
Set-TextValue $ws.Range('E3') '  -0.08%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.74%  '

# Row 5
Set-TextValue $ws.Range('D5') '325.82'
Set-TextValue $ws.Range('E5') '  -0.20%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.4829'
Set-TextValue $ws.Range('E7') '  +1.34%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.4067'
Set-TextValue $ws.Range('E8') '  -0.61%  '

# Row 9
Set-TextValue $ws.Range('E9') '  +1.59%  '

# Row 10
Set-TextValue $ws.Range('D10') '1.012'
Set-TextValue $ws.Range('E10') '  +0.25%  '

# Row 11
Set-TextValue $ws.Range('D11') '23.37'
Set-TextValue $ws.Range('E11') '  +4.27%  '

# Row 12
Set-TextValue $ws.Range('D12') '1.914.96'
Set-TextValue $ws.Range('E12') '  +0.43%  '

# Row 13
Set-TextValue $ws.Range('D13') '6.005'
Set-TextValue $ws.Range('E13') '  +1.23%  '

# Row 14
Set-TextValue $ws.Range('D14') '7.110'
Set-TextValue $ws.Range('E14') '  -0.53%  '

# Row 15
Set-TextValue $ws.Range('D15') '90.42'
Set-TextValue $ws.Range('E15') '  +1.00%  '

# Row 16
Set-TextValue $ws.Range('D16') '0.06809'
Set-TextValue $ws.Range('E16') '  +3.12%  '

# Row 17
Set-TextValue $ws.Range('E17') '  +0.65%  '

# Row 18
Set-TextValue $ws.Range('E18') '  +0.89%  '

# Row 19
Set-TextValue $ws.Range('D19') '17.72'
Set-TextValue $ws.Range('E19') '  -0.04%  '

# Row 20
Set-TextValue $ws.Range('E20') '  +0.66%  '

# Row 21
Set-TextValue $ws.Range('D21') '29.558.20'

# Row 22
Set-TextValue $ws.Range('D22') '5.616'
Set-TextValue $ws.Range('E22') '  +1.37%  '

# Row 23
Set-TextValue $ws.Range('D23') '11.84'
Set-TextValue $ws.Range('E23') '  +2.92%  '

# Row 24
Set-TextValue $ws.Range('E24') '  -1.37%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.151.88'
Set-TextValue $ws.Range('E25') '  +0.83%  '

# Row 26
Set-TextValue $ws.Range('D26') '155.52'
Set-TextValue $ws.Range('E26') '  +0.76%  '

# Row 27
Set-TextValue $ws.Range('B27') 'EthereumClassic'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D27') '20.08'
Set-TextValue $ws.Range('E27') '  +1.39%  '

# Row 28
Set-TextValue $ws.Range('B28') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D28') '6.354'
Set-TextValue $ws.Range('E28') '  +7.78%  '

# Row 29
Set-TextValue $ws.Range('D29') '2.100'
Set-TextValue $ws.Range('E29') '  -1.46%  '

# Row 30
Set-TextValue $ws.Range('D30') '119.54'
Set-TextValue $ws.Range('E30') '  +1.65%  '

# Row 31
Set-TextValue $ws.Range('D31') '1.027'
Set-TextValue $ws.Range('E31') '  -2.12%  '

# Row 32
Set-TextValue $ws.Range('D32') '0.09574'
Set-TextValue $ws.Range('E32') '  +0.48%  '

# Row 33
Set-TextValue $ws.Range('D33') '5.532'
Set-TextValue $ws.Range('E33') '  +2.73%  '

# Row 34
Set-TextValue $ws.Range('D34') '3.558'
Set-TextValue $ws.Range('E34') '  -0.38%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.391'

# Row 36
Set-TextValue $ws.Range('D36') '0.02268'
Set-TextValue $ws.Range('E36') '  +0.55%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.06116'
Set-TextValue $ws.Range('E37') '  -0.05%  '

# Row 38
Set-TextValue $ws.Range('D38') '1.178'
Set-TextValue $ws.Range('E38') '  +0.33%  '

# Row 39
Set-TextValue $ws.Range('B39') 'Aptos'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D39') '10.78'
Set-TextValue $ws.Range('E39') '  +6.26%  '

# Row 40
Set-TextValue $ws.Range('B40') 'TheSandbox'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D40') '0.5934'
Set-TextValue $ws.Range('E40') '  +1.01%  '

# Row 41
Set-TextValue $ws.Range('D41') '7.922'
Set-TextValue $ws.Range('E41') '  -4.83%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.1855'
Set-TextValue $ws.Range('E42') '  +0.74%  '

# Row 43
Set-TextValue $ws.Range('D43') '2.446'
Set-TextValue $ws.Range('E43') '  -3.77%  '

# Row 44
Set-TextValue $ws.Range('D44') '1.286'
Set-TextValue $ws.Range('E44') '  +0.15%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.07732'
Set-TextValue $ws.Range('E45') '  -3.45%  '

# Row 46
Set-TextValue $ws.Range('E46') '  +2.81%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.5561'
Set-TextValue $ws.Range('E47') '  +0.32%  '

# Row 48
Set-TextValue $ws.Range('D48') '1.945'
Set-TextValue $ws.Range('E48') '  +0.85%  '

# Row 49
Set-TextValue $ws.Range('D49') '115.62'
Set-TextValue $ws.Range('E49') '  +2.27%  '

# Row 50
Set-TextValue $ws.Range('D50') '72.66'
Set-TextValue $ws.Range('E50') '  +1.63%  '

# Row 51
Set-TextValue $ws.Range('D51') '1.054'
Set-TextValue $ws.Range('E51') '  +2.09%  '

